# Auto-generated edit script: updates FFXIV Bahamut profit sheets with refreshed market data
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1004.2941
$ws.Range("I28").Value = 1062.3125
$ws.Range("J28").Value = 76
$ws.Range("K28").Value = 1062.3125
$ws.Range("L28").Value = 76
$ws.Range("M28").Value = -577.3125
$ws.Range("N28").Value = -1046
$ws.Range("H40").Value = 2062.5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2062.5
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2062.5
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2412.5
$ws.Range("H64").Value = 5243
$ws.Range("I64").Value = 5140.2
$ws.Range("K64").Value = 5140.2
$ws.Range("M64").Value = -4892.2
$ws.Range("H67").Value = 5243
$ws.Range("I67").Value = 5140.2
$ws.Range("K67").Value = 5140.2
$ws.Range("M67").Value = -4282.2
$ws.Range("H74").Value = 3711444.8
$ws.Range("I74").Value = 4722602
$ws.Range("J74").Value = 3868.3333
$ws.Range("K74").Value = 4722602
$ws.Range("L74").Value = 3868.3333
$ws.Range("M74").Value = -4721666
$ws.Range("N74").Value = -5740.3333
$ws.Range("H76").Value = 60002750
$ws.Range("I76").Value = 62502740
$ws.Range("K76").Value = 62502740
$ws.Range("M76").Value = -62502425
$ws.Range("H77").Value = 3711444.8
$ws.Range("I77").Value = 4722602
$ws.Range("J77").Value = 3868.3333
$ws.Range("K77").Value = 23613010
$ws.Range("L77").Value = 19341.6665
$ws.Range("M77").Value = -23608330
$ws.Range("N77").Value = -28701.6665
$ws.Range("H79").Value = 60002750
$ws.Range("I79").Value = 62502740
$ws.Range("K79").Value = 62502740
$ws.Range("M79").Value = -62501648
$ws.Range("H129").Value = 998.7018
$ws.Range("J129").Value = 1132.2766
$ws.Range("L129").Value = 3396.8298
$ws.Range("N129").Value = -13396.8298
$ws.Range("H138").Value = 3157.5557
$ws.Range("I138").Value = 1118.9722
$ws.Range("J138").Value = 4322.4604
$ws.Range("K138").Value = 3356.9166
$ws.Range("L138").Value = 12967.3812
$ws.Range("M138").Value = 1783.0834
$ws.Range("N138").Value = -23247.3812

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1908.0227
$ws.Range("I132").Value = 1511.9656
$ws.Range("J132").Value = 2673.7334
$ws.Range("K132").Value = 4535.8968
$ws.Range("L132").Value = 8021.2002
$ws.Range("M132").Value = -2005.8968
$ws.Range("N132").Value = -13081.2002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 2100
$ws.Range("I33").Value = 2100
$ws.Range("K33").Value = 2100
$ws.Range("M33").Value = -1764
$ws.Range("H105").Value = 10102
$ws.Range("I105").Value = 10102
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 10102
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -8355
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22120.54
$ws.Range("I31").Value = 1650.9143
$ws.Range("J31").Value = 69883
$ws.Range("K31").Value = 1650.9143
$ws.Range("L31").Value = 69883
$ws.Range("M31").Value = -1355.9143
$ws.Range("N31").Value = -70473
$ws.Range("H34").Value = 22120.54
$ws.Range("I34").Value = 1650.9143
$ws.Range("J34").Value = 69883
$ws.Range("K34").Value = 1650.9143
$ws.Range("L34").Value = 69883
$ws.Range("M34").Value = -1448.9143
$ws.Range("N34").Value = -70287
$ws.Range("H58").Value = 4295.8057
$ws.Range("I58").Value = 1463.7587
$ws.Range("K58").Value = 1463.7587
$ws.Range("M58").Value = -1260.7587
$ws.Range("H105").Value = 3629.8462
$ws.Range("I105").Value = 4124.5454
$ws.Range("K105").Value = 4124.5454
$ws.Range("M105").Value = -2377.5454
$ws.Range("H136").Value = 4295.8057
$ws.Range("I136").Value = 1463.7587
$ws.Range("K136").Value = 4391.2761
$ws.Range("M136").Value = -1841.2761

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 124.25
$ws.Range("I14").Value = 124.25
$ws.Range("K14").Value = 372.75
$ws.Range("M14").Value = -199.75
$ws.Range("H114").Value = 1621.5714
$ws.Range("J114").Value = 1707.8462
$ws.Range("L114").Value = 5123.5386
$ws.Range("N114").Value = -11631.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16245.061
$ws.Range("I70").Value = 18356.357
$ws.Range("K70").Value = 18356.357
$ws.Range("M70").Value = -18086.357
$ws.Range("H73").Value = 16245.061
$ws.Range("I73").Value = 18356.357
$ws.Range("K73").Value = 18356.357
$ws.Range("M73").Value = -17420.357
$ws.Range("H80").Value = 3410.5356
$ws.Range("I80").Value = 3437.2917
$ws.Range("J80").Value = 3250
$ws.Range("K80").Value = 3437.2917
$ws.Range("L80").Value = 3250
$ws.Range("M80").Value = -2439.2917
$ws.Range("N80").Value = -5246
$ws.Range("H83").Value = 3410.5356
$ws.Range("I83").Value = 3437.2917
$ws.Range("J83").Value = 3250
$ws.Range("K83").Value = 17186.4585
$ws.Range("L83").Value = 16250
$ws.Range("M83").Value = -12194.4585
$ws.Range("N83").Value = -26234

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 182.8125
$ws.Range("I55").Value = 182.55556
$ws.Range("J55").Value = 183.14285
$ws.Range("K55").Value = 182.55556
$ws.Range("L55").Value = 183.14285
$ws.Range("M55").Value = -9.555560000000014
$ws.Range("N55").Value = -529.14285
$ws.Range("H82").Value = 5200.6665
$ws.Range("I82").Value = 5200.6665
$ws.Range("K82").Value = 5200.6665
$ws.Range("M82").Value = -4839.6665
$ws.Range("H85").Value = 5200.6665
$ws.Range("I85").Value = 5200.6665
$ws.Range("K85").Value = 5200.6665
$ws.Range("M85").Value = -3952.6665
$ws.Range("H110").Value = 43090.668
$ws.Range("J110").Value = 43090.668
$ws.Range("L110").Value = 43090.668
$ws.Range("N110").Value = -51270.668
$ws.Range("H122").Value = 2518.8572
$ws.Range("I122").Value = 2487.6365
$ws.Range("K122").Value = 7462.9095
$ws.Range("M122").Value = -5012.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1991.8182
$ws.Range("I132").Value = 1753.0667
$ws.Range("J132").Value = 2503.4285
$ws.Range("K132").Value = 5259.2001
$ws.Range("L132").Value = 7510.2855
$ws.Range("M132").Value = -2729.2001
$ws.Range("N132").Value = -12570.2855
